$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "statut" column (A) uses 4 symbolic values that don't render well as
# emoji in some Excel setups. Replace them with safer representations:
#   book-blue   -> warning sign
#   book-red    -> "-3"
#   book-orange -> "+3"
#   book-green  -> check mark
$oldWarn = "📘"
$oldMinus = "📕"
$oldPlus = "📙"
$oldCheck = "📗"

$newWarn = "⚠️"
$newMinus = "-3"
$newPlus = "+3"
$newCheck = "✅"

$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1
$headerRow = $firstRow

for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()

    if ($val -eq $oldWarn) {
        $cell.Value = $newWarn
    }
    elseif ($val -eq $oldCheck) {
        $cell.Value = $newCheck
    }
    elseif ($val -eq $oldMinus) {
        # Force text storage so the leading "-" is kept literally and Excel
        # doesn't reinterpret the content as a numeric value.
        $cell.NumberFormat = "@"
        $cell.Value = $newMinus
    }
    elseif ($val -eq $oldPlus) {
        # Force text storage so the leading "+" is kept literally and Excel
        # doesn't reinterpret the content as a numeric value (and drop the sign).
        $cell.NumberFormat = "@"
        $cell.Value = $newPlus
    }
}
